$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel
# are first switched to Text format so the literal string is preserved.

$ws.Range('D2').Value = '41.369.78'
$ws.Range('E2').Value = '  +3.91%  '
$ws.Range('D3').Value = '2.217.28'
$ws.Range('E3').Value = '  +2.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.99'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.621'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.79'
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.400'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.23'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +5.40%  '
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '2.546.72'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.53'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.54'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '2.223.22'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '41.255.68'
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.56'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').Value = '0.0₃0896'
$ws.Range('E21').Value = '  +5.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.05'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.64'
$ws.Range('E23').Value = '  +9.04%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.38'
$ws.Range('E27').Value = '  -1.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.29'
$ws.Range('E28').Value = '  -2.96%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.89'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  -2.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.59'
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.96'
$ws.Range('E34').Value = '  +5.86%  '
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.54'
$ws.Range('E37').Value = '  -5.48%  '
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.36'
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('B40').Value = 'TerraClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000243'
$ws.Range('E40').Value = '  +26.66%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.996'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('E42').Value = '  -4.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0236'
$ws.Range('E43').Value = '  +3.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.57'
$ws.Range('E44').Value = '  +8.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0979'
$ws.Range('E45').Value = '  +6.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.78'
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').Value = '1.464.63'
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.51'
$ws.Range('E49').Value = '  -6.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.79'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('E51').Value = '  -1.63%  '
